# Auto-generated edit script updating cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.228.75"
$ws.Range("E2").Value = "  +6.02%  "
$ws.Range("D3").Value = "'2.447.94"
$ws.Range("E3").Value = "  +6.34%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'567.35"
$ws.Range("E5").Value = "  +4.95%  "
$ws.Range("D6").Value = "'143.40"
$ws.Range("E6").Value = "  +11.73%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "  +3.72%  "
$ws.Range("D9").Value = "'2.448.86"
$ws.Range("E9").Value = "  +6.44%  "
$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  +5.22%  "
$ws.Range("D11").Value = "'5.77"
$ws.Range("E11").Value = "  +5.16%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "'0.354"
$ws.Range("E13").Value = "  +7.22%  "
$ws.Range("D14").Value = "'26.42"
$ws.Range("E14").Value = "  +14.64%  "
$ws.Range("D15").Value = "'2.877.82"
$ws.Range("E15").Value = "  +6.13%  "
$ws.Range("D16").Value = "'63.104.48"
$ws.Range("E16").Value = "  +6.07%  "
$ws.Range("D17").Value = "'0.0000144"
$ws.Range("E17").Value = "  +9.72%  "
$ws.Range("D18").Value = "'2.441.80"
$ws.Range("E18").Value = "  +5.51%  "
$ws.Range("D19").Value = "'11.28"
$ws.Range("E19").Value = "  +8.62%  "
$ws.Range("D20").Value = "'342.32"
$ws.Range("E20").Value = "  +10.74%  "
$ws.Range("E21").Value = "  +6.75%  "
$ws.Range("D22").Value = "'6.80"
$ws.Range("E22").Value = "  +5.22%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'65.39"
$ws.Range("E24").Value = "  +3.97%  "
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'1.53"
$ws.Range("E27").Value = "  +14.73%  "
$ws.Range("D28").Value = "'8.22"
$ws.Range("E28").Value = "  +6.81%  "
$ws.Range("E29").Value = "  +12.58%  "
$ws.Range("D30").Value = "'6.90"
$ws.Range("E30").Value = "  +19.63%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "'0.0₃0809"
$ws.Range("E31").Value = "  +13.34%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.83"
$ws.Range("E32").Value = "  +8.62%  "
$ws.Range("D33").Value = "'174.67"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = "  +13.29%  "
$ws.Range("E35").Value = "  +6.18%  "
$ws.Range("D36").Value = "'18.77"
$ws.Range("E36").Value = "  +6.29%  "
$ws.Range("D37").Value = "'370.88"
$ws.Range("E37").Value = "  +19.60%  "
$ws.Range("D38").Value = "'4.50"
$ws.Range("E38").Value = "  +12.81%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  +15.22%  "
$ws.Range("D42").Value = "'40.12"
$ws.Range("E42").Value = "  +6.31%  "
$ws.Range("D43").Value = "'150.66"
$ws.Range("E43").Value = "  +11.00%  "
$ws.Range("D44").Value = "'3.74"
$ws.Range("E44").Value = "  +10.11%  "
$ws.Range("D45").Value = "'20.81"
$ws.Range("E45").Value = "  +12.92%  "
$ws.Range("D46").Value = "'0.597"
$ws.Range("E46").Value = "  +5.81%  "
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").Value = "'0.0523"
$ws.Range("E48").Value = "  +7.32%  "
$ws.Range("D49").Value = "'0.0226"
$ws.Range("E49").Value = "  +7.03%  "
$ws.Range("D50").Value = "'18.04"
$ws.Range("E50").Value = "  +8.63%  "
$ws.Range("D51").Value = "'0.0₆0226"
$ws.Range("E51").Value = "  +4.35%  "
